# Apply the "bad_int_value" fixture edit: the documentType column (G) is
# shifted left into F (replacing the old "numberOfPages" column), the
# numeric "page count" sample value moves to a new trailing "numberOfVolumes"
# column (G), and the selection moves from F4 to G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 first, so the brand-new "a" string is appended to the shared
#     string table before "numberOfVolumes" (matches original authoring order).
$ws.Range("F3").Value = "BOOK"
$ws.Range("G3").Value = "a"

# --- Header row: old F1 "numberOfPages" is dropped, documentType moves from
#     G1 to F1, and a new "numberOfVolumes" header is introduced in G1.
$ws.Range("F1").Value = "documentType"
$ws.Range("G1").Value = "numberOfVolumes"

# --- Row 2: documentType value moves from G2 to F2 (as text), and the
#     numeric sample value (previously under "numberOfPages") moves to G2.
$ws.Range("F2").Value = "OTHER"
$ws.Range("G2").Value = 1234

# --- Row 4: documentType value moves from G4 to F4; G4 becomes empty.
$ws.Range("F4").Value = "THESIS"
$ws.Range("G4").ClearContents()

# --- Selection moves from F4 to G1.
[void]$ws.Range("G1").Select()
